$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header text in L1 from "Tình trạng" to "Mô tả lỗi/Hướng dẫn xử lý"
$ws.Range("L1").Value = "Mô tả lỗi/Hướng dẫn xử lý"

# Reset selection to A1 (removes the saved A2 selection in sheetView)
$ws.Range("A1").Select()
